# updated UI for batch operation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the repeated "acma_check" task label to "inter_transmission_merge"
# (shared string used by B2:B3)
$ws.Range("B2:B3").Value = "inter_transmission_merge"

# Restore the app window to its normal (non-maximized) geometry/position,
# matching the refreshed batch-claim UI layout.
$win = $excel.ActiveWindow
$win.WindowState = -4143  # xlNormal
$win.Left = 36720
$win.Top = -1635
$win.Width = 17280
$win.Height = 8970

# Move the active selection to C4, matching the refreshed batch-claim UI state
[void]$ws.Range("C4").Select()
